$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 67, shifting the existing rows 67-69 down to 68-70.
$ws.Rows.Item(67).Insert()

# Populate the newly inserted row 67 with the new weekly price entry.
$ws.Cells.Item(67, 1).Value = 1
$ws.Cells.Item(67, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(67, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(67, 4).Value = 44568
$ws.Cells.Item(67, 5).Value = 15
$ws.Cells.Item(67, 6).Value = "Fruta"
$ws.Cells.Item(67, 7).Value = 100109
$ws.Cells.Item(67, 8).Value = "Uva"
$ws.Cells.Item(67, 9).Value = 100109001
$ws.Cells.Item(67, 10).Value = "Uva"
$ws.Cells.Item(67, 11).Value = "Superior Seedless"
$ws.Cells.Item(67, 12).Value = "Segunda"
$ws.Cells.Item(67, 13).Value = 300
$ws.Cells.Item(67, 14).Value = 19000
$ws.Cells.Item(67, 15).Value = 20000
$ws.Cells.Item(67, 16).Value = 19500
$ws.Cells.Item(67, 17).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(67, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(67, 19).Value = 1083
$ws.Cells.Item(67, 20).Value = 18
